# Updates cryptos list price/volume values (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.036.74"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "3.142.26"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.35"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.23"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.131.16"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("E11").Value = "  +5.10%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.22"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "3.665.80"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "63.852.77"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "3.141.05"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.18"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  +12.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.15"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "80.89"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +9.66%  "
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("E30").Value = "  +7.12%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.62"
$ws.Range("E34").Value = "  +3.95%  "
$ws.Range("D35").Value = "0.0₃0858"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.85"
$ws.Range("E40").Value = "  +6.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.33"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.34"
$ws.Range("E42").Value = "  +7.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("E43").Value = "  +8.87%  "
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("D45").Value = "2.887.36"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.24"
$ws.Range("E46").Value = "  +11.07%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.76"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  +4.13%  "
